{"js": "// Apply the benchSuite docx stat fix-up:\n//  - Rows 1-3 of the (single-column) table become \"0M\" placeholders.\n//  - Ten new single-value rows are inserted right after (old) row 3.\n//  - The three trailing multi-value (tab separated) rows collapse down to\n//    a single summary number each (\"100\", \"0\", \"396\").\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\n// Helper: replace the full text of a table cell's first paragraph while\n// keeping the existing run formatting (rFonts / sz) intact.\nfunction setCellText(cell, text) {\n  const range = cell.body.paragraphs.getFirst().getRange();\n  range.insertText(text, \"Replace\");\n}\n\nconst rows = table.rows.items;\n\n// --- Rows 1-3: \"100\" -> \"0M\", \"0\" -> \"0M\", \"396\" -> \"0M\" ---\nconst row1 = rows[0];\nconst row2 = rows[1];\nconst row3 = rows[2];\nrow1.cells.load(\"items\");\nrow2.cells.load(\"items\");\nrow3.cells.load(\"items\");\nawait context.sync();\n\nsetCellText(row1.cells.items[0], \"0M\");\nsetCellText(row2.cells.items[0], \"0M\");\nsetCellText(row3.cells.items[0], \"0M\");\nawait context.sync();\n\n// --- Insert 10 new rows right after (old) row 3 ---\nconst newRowValues = [\n  [\"25\"],\n  [\"0.00003\"],\n  [\"0.00007\"],\n  [\"0.00004\"],\n  [\"0.00001\"],\n  [\"0.00003\"],\n  [\"0.00004\"],\n  [\"0.00004\"],\n  [\"0.00096\"],\n  [\"100.0\"],\n];\nrow3.insertRows(\"After\", newRowValues.length, newRowValues);\nawait context.sync();\n\n// --- Collapse the three trailing tab-separated rows into single values ---\ntable.rows.load(\"items\");\nawait context.sync();\nconst allRows = table.rows.items;\nconst lastIdx = allRows.length - 1;\n\nconst rowA = allRows[lastIdx - 2]; // was \"12\\t0.00003\\t...\\t100.0\"\nconst rowB = allRows[lastIdx - 1]; // was \"10\\t0.00003\\t...\\t100.0\"\nconst rowC = allRows[lastIdx];     // was \"3\\t0.00003\\t...\\t100.0\"\nrowA.cells.load(\"items\");\nrowB.cells.load(\"items\");\nrowC.cells.load(\"items\");\nawait context.sync();\n\nsetCellText(rowA.cells.items[0], \"100\");\nsetCellText(rowB.cells.items[0], \"0\");\nsetCellText(rowC.cells.items[0], \"396\");\nawait context.sync();\n", "ps1": "# Apply the benchSuite docx stat fix-up:\n#  - Rows 1-3 of the (single-column) table become \"0M\" placeholders.\n#  - Ten new single-value rows are inserted right after (old) row 3.\n#  - The three trailing multi-value (tab separated) rows collapse down to\n#    a single summary number each (\"100\", \"0\", \"396\").\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# --- Rows 1-3: \"100\" -> \"0M\", \"0\" -> \"0M\", \"396\" -> \"0M\" ---\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n\n# --- Insert 10 new rows right after (old) row 3 ---\n# $t.Rows.Add(beforeRow) inserts immediately before \"beforeRow\", so walking\n# the values in reverse and always inserting before the same anchor (the\n# old row 4) lands them in the correct forward order.\n$newValues = @(\"25\", \"0.00003\", \"0.00007\", \"0.00004\", \"0.00001\", \"0.00003\", \"0.00004\", \"0.00004\", \"0.00096\", \"100.0\")\n\n$afterRow = $t.Rows(4)\nfor ($i = $newValues.Count - 1; $i -ge 0; $i--) {\n    $newRow = $t.Rows.Add($afterRow)\n    $newRow.Cells(1).Range.Text = $newValues[$i]\n}\n\n# --- Collapse the three trailing tab-separated rows into single values ---\n$rowCount = $t.Rows.Count\n$t.Cell($rowCount - 2, 1).Range.Text = \"100\"\n$t.Cell($rowCount - 1, 1).Range.Text = \"0\"\n$t.Cell($rowCount, 1).Range.Text = \"396\"\n"}
